# Typocrypha spellDictionary.xlsx edit
# Commit: "Added ally damage, stun/stagger, gauge mechanics, and UI"
# Relevant spreadsheet change: add a new "teamcrit" debug spell row
# ("Added teamcrit debug spell that auto-crits both your allies")
# as a new row 12 on the spellDictionary sheet, pushing all subsequent
# rows (old 12-31) down by one (to 13-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 12 (shifts rows 12+ down by one, and
# Excel automatically grows the used range / dimension for us).
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the "teamcrit" debug spell.
# Columns: A=name, B=spell type, C=spellbook description, D=power,
# E=cooldown, F=hit %, G=crit %, H=status%, I=pattern, J=buff/debuffs, K=buff%
$ws.Range("A12").Value = "teamcrit"
$ws.Range("B12").Value = "attack"
$ws.Range("C12").Value = "crits own team (debug)"
$ws.Range("D12").Value = 25
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 100
$ws.Range("G12").Value = 100
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = "lr"
$ws.Range("J12").Value = "None"
$ws.Range("K12").Value = 0

# Move/restore the sheet's active cell selection to K12 (matches the
# author's cursor position after adding the new row).
$ws.Activate()
$ws.Range("K12").Select()
